$wb = $excel.ActiveWorkbook

# Sheets involved
$wsAbout = $wb.Worksheets.Item("About")
$wsBts   = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$wsPax   = $wb.Worksheets.Item("AVLo-passengers")
$wsFrt   = $wb.Worksheets.Item("AVLo-freight")

# Remember the current (pre-edit) active selection on "About" so we can
# restore it once we're done poking at the other sheets (Select()/Activate()
# necessarily move the active sheet/cell).
$aboutSelection = $wsAbout.Application.ActiveCell

# --- Insert a new row 37 into "BTS NTS Modal Profile Data" --------------
# This is a true row insert: everything at/after row 37 shifts down by one,
# and Excel auto-adjusts in-sheet + cross-sheet formula references that
# point at the shifted rows.
$wsBts.Rows.Item(37).Insert()

# New row 37 content: a per-railcar adjustment of the row-36 weighted
# vehicle-loading value (divide by 10 cars/locomotive).
$wsBts.Range("A37").Value = "weighted value, adjusted for number of train cars per locomotive"
$wsBts.Range("A37").WrapText = $true

$wsBts.Range("B37").Formula = "=B36/10"
$wsBts.Range("B37").NumberFormat = "0"

# --- AVLo-passengers: point the vehicle-loading lookup at the new --------
# per-car-adjusted row instead of the raw per-locomotive row. (The other
# reference on this sheet, to the old row 59/now row 60, already shifted
# automatically above.)
$wsPax.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B37"

# --- View state -----------------------------------------------------------
# Reflect the workbook being left with the BTS sheet scrolled to/selecting
# B38, and the freight sheet selecting B6.
$wsBts.Activate()
$wsBts.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 7

$wsFrt.Activate()
$wsFrt.Range("B6").Select()

# Restore "About" as the active sheet/selection, matching the unchanged
# pre-edit state there.
$wsAbout.Activate()
$aboutSelection.Select()
